$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B-column values (rows 2-148) per the diff
$ws.Cells.Item(2, 2).Value = 1.098540744384312
$ws.Cells.Item(3, 2).Value = 1.200129743756086
$ws.Cells.Item(4, 2).Value = 1.025641025641024
$ws.Cells.Item(5, 2).Value = 0.4124365482233533
$ws.Cells.Item(6, 2).Value = 0.9004739336492815
$ws.Cells.Item(7, 2).Value = 1.174260216063888
$ws.Cells.Item(8, 2).Value = 2.073661405137726
$ws.Cells.Item(9, 2).Value = 0.4851425106124861
$ws.Cells.Item(10, 2).Value = 2.112251056125549
$ws.Cells.Item(11, 2).Value = 1.773049645390058
$ws.Cells.Item(12, 2).Value = 2.758420441347268
$ws.Cells.Item(13, 2).Value = -0.4801581697500401
$ws.Cells.Item(14, 2).Value = -0.156094792110117
$ws.Cells.Item(15, 2).Value = 1.43547470153497
$ws.Cells.Item(16, 2).Value = 1.373125963289883
$ws.Cells.Item(17, 2).Value = -0.6219765031098669
$ws.Cells.Item(18, 2).Value = -0.2364394993045921
$ws.Cells.Item(19, 2).Value = -0.2788233653980243
$ws.Cells.Item(20, 2).Value = -0.7968684468055465
$ws.Cells.Item(21, 2).Value = 0.01409244644871071
$ws.Cells.Item(22, 2).Value = 0.5917993518388075
$ws.Cells.Item(23, 2).Value = -0.08404538450763731
$ws.Cells.Item(24, 2).Value = 1.415954016542836
$ws.Cells.Item(25, 2).Value = 0.5529444290848651
$ws.Cells.Item(26, 2).Value = 0.5911465493538725
$ws.Cells.Item(27, 2).Value = 1.107011070110704
$ws.Cells.Item(28, 2).Value = -0.3919978372533202
$ws.Cells.Item(29, 2).Value = 0.8006513773917809
$ws.Cells.Item(30, 2).Value = 0.2019386106623663
$ws.Cells.Item(31, 2).Value = 0.04030632809349311
$ws.Cells.Item(32, 2).Value = -0.7117915659414276
$ws.Cells.Item(33, 2).Value = 1.366157175706737
$ws.Cells.Item(34, 2).Value = 0.3869762476648069
$ws.Cells.Item(35, 2).Value = 0.9171872923035993
$ws.Cells.Item(36, 2).Value = -0.513698630136987
$ws.Cells.Item(37, 2).Value = 1.11214087117702
$ws.Cells.Item(38, 2).Value = 0.3928244074898483
$ws.Cells.Item(39, 2).Value = 0.7434459371331592
$ws.Cells.Item(40, 2).Value = 0.95805282237184
$ws.Cells.Item(41, 2).Value = -0.4616568350859187
$ws.Cells.Item(42, 2).Value = 0.4895645452202981
$ws.Cells.Item(43, 2).Value = -0.02564102564102054
$ws.Cells.Item(44, 2).Value = 1.154142087714788
$ws.Cells.Item(45, 2).Value = -0.1267748478701753
$ws.Cells.Item(46, 2).Value = 1.370906321401369
$ws.Cells.Item(47, 2).Value = 0.4758327072376602
$ws.Cells.Item(48, 2).Value = 1.545363908275186
$ws.Cells.Item(49, 2).Value = 0.7977417771232099
$ws.Cells.Item(50, 2).Value = 0.07305491294289818
$ws.Cells.Item(51, 2).Value = -0.4258425599221248
$ws.Cells.Item(52, 2).Value = 2.126099706744862
$ws.Cells.Item(53, 2).Value = -0.2512562814070277
$ws.Cells.Item(54, 2).Value = 0
$ws.Cells.Item(55, 2).Value = -0.3118627803766404
$ws.Cells.Item(56, 2).Value = -0.4692576103958616
$ws.Cells.Item(57, 2).Value = 0.3868471953578426
$ws.Cells.Item(58, 2).Value = 0.68641618497109
$ws.Cells.Item(59, 2).Value = -0.2750867121157804
$ws.Cells.Item(60, 2).Value = -1.355241065003592
$ws.Cells.Item(61, 2).Value = 0.1458966565349599
$ws.Cells.Item(62, 2).Value = 0.9105256768240864
$ws.Cells.Item(63, 2).Value = 0.1203079884504263
$ws.Cells.Item(64, 2).Value = -0.1562124489305401
$ws.Cells.Item(65, 2).Value = 0.589721988205554
$ws.Cells.Item(66, 2).Value = -0.323043790380469
$ws.Cells.Item(67, 2).Value = -0.03601008282319185
$ws.Cells.Item(68, 2).Value = 0.1560999039385152
$ws.Cells.Item(69, 2).Value = 0.5994485073732168
$ws.Cells.Item(70, 2).Value = 0.7865570253843364
$ws.Cells.Item(71, 2).Value = 0.3783847700130157
$ws.Cells.Item(72, 2).Value = 1.119095299799744
$ws.Cells.Item(73, 2).Value = 1.712488350419384
$ws.Cells.Item(74, 2).Value = 0.7559271561104073
$ws.Cells.Item(75, 2).Value = 1.466409003069235
$ws.Cells.Item(76, 2).Value = 0.1344387183508742
$ws.Cells.Item(77, 2).Value = 0.6601029313045462
$ws.Cells.Item(78, 2).Value = 0.5112815382905501
$ws.Cells.Item(79, 2).Value = 0.7409045670684369
$ws.Cells.Item(80, 2).Value = 0.5817782656421527
$ws.Cells.Item(81, 2).Value = -0.3601440576230474
$ws.Cells.Item(82, 2).Value = -0.547645125958379
$ws.Cells.Item(83, 2).Value = -1.552863436123344
$ws.Cells.Item(84, 2).Value = -4.687325204161537
$ws.Cells.Item(85, 2).Value = 0.2464788732394293
$ws.Cells.Item(86, 2).Value = 0.6205362369745945
$ws.Cells.Item(87, 2).Value = 0.7679776588317392
$ws.Cells.Item(88, 2).Value = 0.7274826789838449
$ws.Cells.Item(89, 2).Value = 2.166685773243151
$ws.Cells.Item(95, 2).Value = -0.01075731497417266
$ws.Cells.Item(96, 2).Value = 0.2151694459386798
$ws.Cells.Item(97, 2).Value = 0.1073537305421302
$ws.Cells.Item(98, 2).Value = 0.2037533512064319
$ws.Cells.Item(99, 2).Value = -0.2996575342465766
$ws.Cells.Item(100, 2).Value = -0.5581794761700257
$ws.Cells.Item(101, 2).Value = 1.165803108808288
$ws.Cells.Item(102, 2).Value = 0.5121638924455868
$ws.Cells.Item(103, 2).Value = 0.1486199575371556
$ws.Cells.Item(104, 2).Value = 1.049395802416785
$ws.Cells.Item(105, 2).Value = -0.02097975453686774
$ws.Cells.Item(106, 2).Value = 0.5560801594795941
$ws.Cells.Item(107, 2).Value = 0.7616861435726103
$ws.Cells.Item(108, 2).Value = -0.2174588381484869
$ws.Cells.Item(109, 2).Value = 0.5188875051888751
$ws.Cells.Item(110, 2).Value = 0.5368573198430684
$ws.Cells.Item(111, 2).Value = 0.5134524543027316
$ws.Cells.Item(112, 2).Value = 0.8786268900694724
$ws.Cells.Item(113, 2).Value = 0.2228073728985202
$ws.Cells.Item(114, 2).Value = 0.3536782538399439
$ws.Cells.Item(115, 2).Value = 0.4631960527640658
$ws.Cells.Item(116, 2).Value = 1.292973839831619
$ws.Cells.Item(117, 2).Value = 0.7124480506629713
$ws.Cells.Item(118, 2).Value = 0.7958341520927513
$ws.Cells.Item(119, 2).Value = 1.01374402963251
$ws.Cells.Item(120, 2).Value = -0.5500337740036604
$ws.Cells.Item(121, 2).Value = 0.7859499320784032
$ws.Cells.Item(122, 2).Value = -0.654664484451725
$ws.Cells.Item(123, 2).Value = 0.4457796298090977
$ws.Cells.Item(124, 2).Value = 0.6560540279787676
$ws.Cells.Item(125, 2).Value = 0.05750982459503716
$ws.Cells.Item(126, 2).Value = 0.3831784653702516
$ws.Cells.Item(127, 2).Value = -0.3340013360053521
$ws.Cells.Item(128, 2).Value = -2.029873611643053
$ws.Cells.Item(129, 2).Value = -8.874120406567631
$ws.Cells.Item(130, 2).Value = 8.67653367653368
$ws.Cells.Item(131, 2).Value = 0.9671370768775328
$ws.Cells.Item(132, 2).Value = -0.6157755840093788
$ws.Cells.Item(133, 2).Value = 2.350511408339877
$ws.Cells.Item(134, 2).Value = 0.08648025367541406
$ws.Cells.Item(135, 2).Value = 0.5376344086021527
$ws.Cells.Item(136, 2).Value = 0.6779984721161268
$ws.Cells.Item(137, 2).Value = 0.1517594612539093
$ws.Cells.Item(138, 2).Value = 0.2935884079931833
$ws.Cells.Item(139, 2).Value = -0.3493862134088806
$ws.Cells.Item(140, 2).Value = -0.473798919738463
$ws.Cells.Item(141, 2).Value = -0.07616871370084576
$ws.Cells.Item(142, 2).Value = 0
$ws.Cells.Item(143, 2).Value = -0.2763220581229216
$ws.Cells.Item(144, 2).Value = -0.1051022358111976
$ws.Cells.Item(145, 2).Value = -0.2582496413199388
$ws.Cells.Item(146, 2).Value = 0.01917913310317992
$ws.Cells.Item(147, 2).Value = 0.1821668264621263
$ws.Cells.Item(148, 2).Value = 0.3062494018566441

# Add new row 149 (copy formatting from row 148, then set values)
$ws.Range("A148").Copy() | Out-Null
$ws.Range("A149").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(149, 1).Value = 45748
$ws.Cells.Item(149, 2).Value = -0.2766911554241067
$excel.CutCopyMode = 0
